$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 53, shifting all rows from
# 53 downward to 55 onward.
$ws.Rows("53:54").Insert()

# --- New row 53 ---
$ws.Range("A53").Value = 7
$ws.Range("B53").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C53").Value = "Ñuble"
$ws.Range("D53").Value = 44924
$ws.Range("E53").Value = 16
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100103
$ws.Range("H53").Value = "Frutos de hueso (carozo)"
$ws.Range("I53").Value = 100103002
$ws.Range("J53").Value = "Ciruela"
$ws.Range("K53").Value = "Black Amber"
$ws.Range("L53").Value = "Especial"
$ws.Range("M53").Value = 60
$ws.Range("N53").Value = 17000
$ws.Range("O53").Value = 17000
$ws.Range("P53").Value = 17000
$ws.Range("Q53").Value = "$/bandeja 18 kilos granel"
$ws.Range("R53").Value = "Provincia de Curicó"
$ws.Range("S53").Value = 944
$ws.Range("T53").Value = 18

# --- New row 54 ---
$ws.Range("A54").Value = 7
$ws.Range("B54").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C54").Value = "Ñuble"
$ws.Range("D54").Value = 44924
$ws.Range("E54").Value = 16
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100103
$ws.Range("H54").Value = "Frutos de hueso (carozo)"
$ws.Range("I54").Value = 100103002
$ws.Range("J54").Value = "Ciruela"
$ws.Range("K54").Value = "Black Amber"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 120
$ws.Range("N54").Value = 15000
$ws.Range("O54").Value = 16000
$ws.Range("P54").Value = 15500
$ws.Range("Q54").Value = "$/bandeja 18 kilos granel"
$ws.Range("R54").Value = "Provincia de Curicó"
$ws.Range("S54").Value = 861
$ws.Range("T54").Value = 18
